# Applies the "semana epidemiologica 43 de 2025" update to poisson.xlsx
# Updates Esperado (C), Observado (D) and valor p (E) columns for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value
$changes = @{
    3  = @{ D = 6 }
    5  = @{ C = 5;  D = 11; E = 0.01 }
    6  = @{ D = 2;  E = 0.27 }
    7  = @{ C = 3;  D = 8;  E = 0.01 }
    9  = @{ C = 38; D = 50; E = 0.01 }
    11 = @{ D = 2;  E = 0.18 }
    12 = @{ C = 9;  D = 3;  E = 0.01 }
    16 = @{ C = 0;  E = 1 }
    17 = @{ C = 9;  D = 19; E = 0 }
    18 = @{ C = 2;  E = 0.14 }
    19 = @{ C = 8;  D = 2;  E = 0.01 }
    20 = @{ C = 0;  E = 1 }
    22 = @{ C = 1;  D = 3;  E = 0.06 }
    23 = @{ D = 0;  E = 1 }
    25 = @{ C = 7;  D = 10; E = 0.07000000000000001 }
    29 = @{ C = 1;  E = 0.37 }
    30 = @{ D = 0;  E = 1 }
    31 = @{ D = 1 }
    33 = @{ C = 7;  D = 3;  E = 0.05 }
    34 = @{ C = 10; D = 1;  E = 0 }
    35 = @{ C = 8;  D = 9;  E = 0.12 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $ws.Range($cellRef).Value = $cols[$col]
    }
}
